# Slide 2, shape "副標題 2" (the subtitle placeholder) holds the bullet list.
# The 11th paragraph reads: "How to create a cookie? setcookie()." and is
# split across three runs: "How to create a cookie? s" | "etcookie" | "()."
# Fix the typo split so the word "setcookie" is spelled correctly and sits
# entirely in the second run.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$para = $tr.Paragraphs(11, 1)

$run1 = $para.Runs(1, 1)
$run1.Text = "How to create a cookie? "

$run2 = $para.Runs(2, 1)
$run2.Text = "setcookie"
